# Petty Cash Book 2021 - 8 Apr 2021 midday update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19: additional 260000 debit on Wages Expense
$ws.Range("D19").Formula = "=60000+260000"

# Row 21: TRANSFER BCA - additional amounts debited
$ws.Range("D21").Formula = "=504000+1057500+750000+3040000+40000000+35111000+11215000+3324000+2600000"

# Row 22: new A/R entry (credit)
$ws.Range("B22").Value = "A/R"
$ws.Range("C22").Formula = "=40000000+855000+11215000+2600000+46251500"

# Row 23: new SALES - cash/retail entry (credit)
$ws.Range("B23").Value = "SALES - cash/retail"
$ws.Range("C23").Formula = "=11216525+42768975-46251500"

# Row 24: new SELISIH - lebih entry (credit)
$ws.Range("B24").Value = "SELISIH - lebih"
$ws.Range("C24").Value = 30000

# Row 25: new SETOR KE BANK entry (debit)
$ws.Range("B25").Value = "SETOR KE BANK"
$ws.Range("D25").Value = 10000000

# Row 26: new dated entry 7 Apr 2021, Wages Expense
$ws.Range("A26").Value = 44293
$ws.Range("B26").Value = "Wages Expense"
$ws.Range("D26").Formula = "=60000+240000"

# Row 27: new TRANSFER BCA entry (debit)
$ws.Range("B27").Value = "TRANSFER BCA"
$ws.Range("D27").Formula = "=300000+13000000+8311000"

# Row 28: new A/R entry (credit)
$ws.Range("B28").Value = "A/R"
$ws.Range("C28").Formula = "=13000000+8751000"

# Row 29: new A/P entry (debit)
$ws.Range("B29").Value = "A/P"
$ws.Range("D29").Formula = "=478000"

# Row 30: new SALES - cash/retail entry (credit)
$ws.Range("B30").Value = "SALES - cash/retail"
$ws.Range("C30").Formula = "=8162475+4973525-8751000"

# Row 31: new SETOR KE BANK entry (debit)
$ws.Range("B31").Value = "SETOR KE BANK"
$ws.Range("D31").Formula = "=4000000"

# Row 32: new dated entry 8 Apr 2021
$ws.Range("A32").Value = 44294

# Update selection to reflect the latest edited cell
$ws.Range("B32").Select()
